# Update cryptocurrency price/volume data to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text updates (values that Excel will not mistake for numbers).
$ws.Range("D2").Value = '26.280.12'
$ws.Range("E2").Value = '  +0.32%  '
$ws.Range("D3").Value = '1.588.53'
$ws.Range("E3").Value = '  +0.43%  '
$ws.Range("E4").Value = '  -0.20%  '
$ws.Range("E5").Value = '  +1.18%  '
$ws.Range("E6").Value = '  +1.58%  '
$ws.Range("E7").Value = '  -0.19%  '
$ws.Range("E8").Value = '  +0.16%  '
$ws.Range("E9").Value = '  -0.24%  '
$ws.Range("D12").Value = '1.812.66'
$ws.Range("E12").Value = '  +0.48%  '
$ws.Range("D13").Value = '1.582.70'
$ws.Range("E13").Value = '  -2.36%  '
$ws.Range("E14").Value = '  -0.45%  '
$ws.Range("E15").Value = '  +1.05%  '
$ws.Range("E16").Value = '  -0.11%  '
$ws.Range("D17").Value = '26.292.66'
$ws.Range("E17").Value = '  +0.38%  '
$ws.Range("D18").Value = '0.0₃0730'
$ws.Range("E18").Value = '  -0.50%  '
$ws.Range("E19").Value = '  +3.04%  '
$ws.Range("E20").Value = '  +3.17%  '
$ws.Range("E21").Value = '  -0.18%  '
$ws.Range("E22").Value = '  +0.85%  '
$ws.Range("E23").Value = '  +1.07%  '
$ws.Range("E24").Value = '  -2.07%  '
$ws.Range("E25").Value = '  +0.20%  '
$ws.Range("E26").Value = '  -0.15%  '
$ws.Range("E27").Value = '  +0.43%  '
$ws.Range("E28").Value = '  -0.26%  '
$ws.Range("E29").Value = '  -0.03%  '
$ws.Range("E30").Value = '  -0.76%  '
$ws.Range("E31").Value = '  +1.05%  '
$ws.Range("E32").Value = '  -0.21%  '
$ws.Range("E33").Value = '  +0.94%  '
$ws.Range("D34").Value = '1.339.95'
$ws.Range("E34").Value = '  +4.74%  '
$ws.Range("E35").Value = '  -1.02%  '
$ws.Range("B36").Value = 'LidoDAOToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("E36").Value = '  -0.14%  '
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("E37").Value = '  -1.62%  '
$ws.Range("E38").Value = '  +0.04%  '
$ws.Range("E39").Value = '  -9.96%  '
$ws.Range("E40").Value = '  +0.23%  '
$ws.Range("E41").Value = '  +3.70%  '
$ws.Range("E42").Value = '  -0.22%  '
$ws.Range("E43").Value = '  +0.35%  '
$ws.Range("E44").Value = '  -0.34%  '
$ws.Range("D45").Value = '1.724.12'
$ws.Range("E45").Value = '  +0.28%  '
$ws.Range("E46").Value = '  -0.76%  '
$ws.Range("E47").Value = '  -0.83%  '
$ws.Range("E48").Value = '  -4.04%  '
$ws.Range("E49").Value = '  -2.61%  '
$ws.Range("E50").Value = '  -0.66%  '
$ws.Range("E51").Value = '  -0.44%  '

# Numeric-looking price strings: force text format first so Excel keeps
# them as literal strings (matching the source data) instead of coercing
# them into floating-point numbers, then restore the default formatting.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.01'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.34'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.42'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.47'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '213.74'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.97'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.15'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.08'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.04'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.19'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.16'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.48'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.598'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '61.77'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '88.16'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
